# SkillRef.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to be called "Property1" is renamed to "DataNode",
# and the author's last selection on that sheet moved to C38 (still in
# the frozen "bottomLeft" pane below the header rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "DataNode"
$ws.Range("C38").Select()
